$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "figures manuscript / absolute n" column (D) values for the rows that
# previously had no entry there. Most rows are "yes"; two are flagged
# differently ("no?" on row 31, "no" on row 36).
$ws.Range("D25").Value = "yes"

$ws.Range("D30").Value = "yes"
$ws.Range("D31").Value = "no?"
$ws.Range("D32").Value = "yes"
$ws.Range("D33").Value = "yes"
$ws.Range("D34").Value = "yes"
$ws.Range("D35").Value = "yes"
$ws.Range("D36").Value = "no"
$ws.Range("D37").Value = "yes"
$ws.Range("D38").Value = "yes"
$ws.Range("D39").Value = "yes"
$ws.Range("D40").Value = "yes"
$ws.Range("D41").Value = "yes"
$ws.Range("D42").Value = "yes"
$ws.Range("D43").Value = "yes"
$ws.Range("D44").Value = "yes"
$ws.Range("D45").Value = "yes"
$ws.Range("D46").Value = "yes"
$ws.Range("D47").Value = "yes"
$ws.Range("D48").Value = "yes"
$ws.Range("D49").Value = "yes"

# Restore the active view: the user had scrolled down and was looking at
# E48 when the workbook was last saved.
$ws.Activate()
[void]$ws.Range("E48").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 24
$window.ScrollColumn = 1
